# LOM3003.xlsx update
# - Removes the standalone row that held only the professor's name
#   (B13/C13 "984972 - Hugo Ricardo Zschommler Sandim"), shifting all
#   rows below it up by one.
# - Re-populates several B/C cells with the (reshuffled) content that the
#   published sheet now shows under each label in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old standalone "984972 - Hugo Ricardo Zschommler Sandim" row.
# Everything from row 14 downward shifts up to become row 13 downward,
# carrying its row heights and styles along with it.
$ws.Rows.Item(13).Delete()

# After the shift, write the new B/C text for the affected rows.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2016"
$ws.Range("C15").Value = "01/01/2016"

$ws.Range("B18").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C18").Value = "984972 - Hugo Ricardo Zschommler Sandim"

$ws.Range("B19").Value = "O aluno será avaliado ao longo do semestre por duas avaliações escritas (P1 e P2) e com pesos iguais."
$ws.Range("C19").Value = "O aluno será avaliado ao longo do semestre por duas avaliações escritas (P1 e P2) e com pesos iguais."

$ws.Range("B20").Value = "Nota Final NF = [P1 + P2]/2"
$ws.Range("C20").Value = "Nota Final NF = [P1 + P2]/2"

$ws.Range("B21").Value = "Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2"
$ws.Range("C21").Value = "Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2"
